$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the row above (row 22) onto the new row (23)
# so the new date/time cells pick up the same cell styles (s="1"/s="3")
# instead of creating brand-new style entries.
$ws.Range("A22:D22").Copy()
$ws.Range("A23").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row's data
$ws.Range("A23").Value2 = 43752
$ws.Range("B23").Value = "AdobeXD GUI Modell fertig gestellt + Vue Navigation Bar informiert"
$ws.Range("C23").Value2 = 0.85416666666666663
$ws.Range("D23").Value2 = 0.88194444444444453

$ws.Range("G12").Select()
